$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing entry (row 49): "4.iterace" -> "4. iterace"
$ws.Range("A49").Value = "kontrola a oponentura 4. iterace ŠIP"

# Append a new activity row
$ws.Range("A50").Value = "aktualizace dokumentace - oprava testovacích scénářů"
$ws.Range("B50").Value = 0.5

# Keep column A sized to fit the widened content
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update selection to mirror the new last empty row
$ws.Range("A51").Select()
